# Added January unemployment data (new column C) to the "Unemployment"
# sheet, plus the small selection-only change on "Small Businesses by
# State" that accompanied it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Small Businesses by State": no data change, only the active cell /
# selection moved to C1.
# ---------------------------------------------------------------------
$wsSmallBiz = $wb.Worksheets.Item("Small Businesses by State")
$wsSmallBiz.Activate() | Out-Null
$wsSmallBiz.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------
# "Unemployment": add a "January" column next to the existing "March"
# column.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Unemployment")

# A style-only donor cell (Arial 14, not bold == style index 7 in the
# original file) so the new blank D cells / trailing row pick up the
# same formatting Excel would have applied automatically.
$blankStyleSrc = $wsSmallBiz.Cells.Item(2, 1)

# Header: C1 = "January", formatted like the existing A1/B1 headers.
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)
$ws.Cells.Item(1, 3).Value = "January"

# D1 stays blank but still carries formatting.
$blankStyleSrc.Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)

# January unemployment rates, state rows 2-52 (same order as column B).
$januaryValues = @(
    6, 2.7, 4.5, 3.5, 3.9, 2.5, 3.7, 4, 5.2, 2.8,
    3.1, 2.7, 2.8, 3.5, 3.1, 2.8, 3.1, 4.3, 5.0999999999999996, 3.1,
    3.3, 2.8, 3.8, 3.2, 5.5, 3.5, 3.5, 3.9, 3.6, 2.6,
    3.8, 4.8, 3.8, 3.6, 2.2999999999999998, 4.0999999999999996, 3.3, 3.3, 4.7, 3.1,
    2.4, 3.4, 3.3, 4.5, 2.5, 2.4, 2.7, 3.9, 5, 3.5,
    3.7
)

for ($i = 0; $i -lt $januaryValues.Count; $i++) {
    $row = 2 + $i

    # Give C<row> the same look as B<row> (number style), then overwrite
    # with the real January value.
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value = $januaryValues[$i]

    # D<row> stays blank but formatted.
    $blankStyleSrc.Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
}

# Trailing blank row 53 under the table (C53:D53), formatted only.
$blankStyleSrc.Copy()
$ws.Cells.Item(53, 3).PasteSpecial(-4122)
$blankStyleSrc.Copy()
$ws.Cells.Item(53, 4).PasteSpecial(-4122)

# Widen the new column to fit the "January" header / values.
$ws.Columns.Item(3).ColumnWidth = 19.330729166666668

# Restore the active sheet/selection to "Unemployment".
$ws.Activate() | Out-Null
$ws.Range("G14").Select() | Out-Null
